# "Use '|' symbol for lane divider"
#
# The "Changes" sheet row for the "Lane" field (row 5) had its Final Value
# updated from the bare status name "Review" to a pipe-delimited
# ART/VS|Art/Status path: "Implementation per ART/VS|Art 1|In Progress".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Changes")
$ws.Range("F5").Value = "Implementation per ART/VS|Art 1|In Progress"

# The author was on the "Changes" sheet when the workbook was saved, so make
# it the active sheet/tab (this also moves tabSelected in the sheetView and
# updates workbookView's activeTab).
$ws.Activate()
